# factor out resolution_z_unit (#817)
#
# 1. Comment on column V (resolution_z_unit) header loses the "(um)" suffix
#    now that "mm" is also a valid unit.
# 2. The "resolution_z_unit list" sheet (used for the column V dropdown)
#    gains a new "mm" option, inserted ahead of the existing "nm"/"um"
#    entries: A1=mm, A2=um, A3=nm.
# 3. The data validation on column V is widened to pull from the 3-row
#    list and the error message is updated to mention all three units.

$wb = $excel.ActiveWorkbook

$wsMain = $wb.Worksheets.Item("Export as TSV")
$wsZUnit = $wb.Worksheets.Item("resolution_z_unit list")

# --- 1. Update the V1 cell comment text -------------------------------
$null = $wsMain.Range("V1").Comment.Text("The unit of incremental distance between image slices.")

# --- 2. Reorder / extend the resolution_z_unit list sheet --------------
# Original: A1=nm, A2=um
# New:      A1=mm, A2=um, A3=nm
$wsZUnit.Range("A2").Value = "um"
$wsZUnit.Range("A3").Value = "nm"
$wsZUnit.Range("A1").Value = "mm"

# --- 3. Widen the data validation list + update its error message ------
$validation = $wsMain.Range("V2:V1048576").Validation
$validation.Modify(3, 1, 1, '''resolution_z_unit list''!$A$1:$A$3')
$validation.ErrorTitle = "Value must come from list"
$validation.ErrorMessage = "Value must be one of: mm / um / nm."
